$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 709.3
$ws.Range("I15").Value = 709.3
$ws.Range("K15").Value = 2127.9
$ws.Range("M15").Value = -1958.9

# Row 28
$ws.Range("H28").Value = 709.9
$ws.Range("I28").Value = 525.7727
$ws.Range("J28").Value = 1216.25
$ws.Range("K28").Value = 525.7727
$ws.Range("L28").Value = 1216.25
$ws.Range("M28").Value = -40.77269999999999
$ws.Range("N28").Value = -2186.25

# Row 86
$ws.Range("H86").Value = 5327.609
$ws.Range("I86").Value = 3045.3333
$ws.Range("K86").Value = 3045.3333
$ws.Range("M86").Value = -1922.3333

# Row 89
$ws.Range("H89").Value = 5327.609
$ws.Range("I89").Value = 3045.3333
$ws.Range("K89").Value = 15226.6665
$ws.Range("M89").Value = -9610.666499999999

# Row 107
$ws.Range("H107").Value = 721.36365
$ws.Range("I107").Value = 393.83334
$ws.Range("J107").Value = 1114.4
$ws.Range("K107").Value = 393.83334
$ws.Range("L107").Value = 1114.4
$ws.Range("M107").Value = 1526.16666
$ws.Range("N107").Value = -4954.4

# Row 116
$ws.Range("H116").Value = 92273.414
$ws.Range("I116").Value = 135985.12
$ws.Range("J116").Value = 4850
$ws.Range("K116").Value = 135985.12
$ws.Range("L116").Value = 4850
$ws.Range("M116").Value = -132543.12
$ws.Range("N116").Value = -11734

# Row 139
$ws.Range("H139").Value = 35000
$ws.Range("I139").Value = 35000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 35000
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -29860

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 4398.222
$ws.Range("I26").Value = 2014
$ws.Range("J26").Value = 9166.666999999999
$ws.Range("K26").Value = 2014
$ws.Range("L26").Value = 9166.666999999999
$ws.Range("M26").Value = -1684
$ws.Range("N26").Value = -9826.666999999999

# Row 74
$ws.Range("H74").Value = 1603.8235
$ws.Range("I74").Value = 702.9737
$ws.Range("J74").Value = 4237.077
$ws.Range("K74").Value = 702.9737
$ws.Range("L74").Value = 4237.077
$ws.Range("M74").Value = 171.0263
$ws.Range("N74").Value = -5985.077

# Row 77
$ws.Range("H77").Value = 1603.8235
$ws.Range("I77").Value = 702.9737
$ws.Range("J77").Value = 4237.077
$ws.Range("K77").Value = 3514.8685
$ws.Range("L77").Value = 21185.385
$ws.Range("M77").Value = 853.1315
$ws.Range("N77").Value = -29921.385

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 25000
$ws.Range("I26").Value = 20000
$ws.Range("K26").Value = 20000
$ws.Range("M26").Value = -19708

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 212.52942
$ws.Range("I107").Value = 157.91667
$ws.Range("J107").Value = 343.6
$ws.Range("K107").Value = 157.91667
$ws.Range("L107").Value = 343.6
$ws.Range("M107").Value = 1762.08333
$ws.Range("N107").Value = -4183.6

# Row 132
$ws.Range("H132").Value = 2533.0344
$ws.Range("I132").Value = 1678.125
$ws.Range("J132").Value = 3585.2307
$ws.Range("K132").Value = 5034.375
$ws.Range("L132").Value = 10755.6921
$ws.Range("M132").Value = -2504.375
$ws.Range("N132").Value = -15815.6921

$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 3236.2666
$ws.Range("I64").Value = 1666
$ws.Range("J64").Value = 3477.8462
$ws.Range("K64").Value = 4998
$ws.Range("L64").Value = 10433.5386
$ws.Range("M64").Value = -4728
$ws.Range("N64").Value = -10973.5386

# Row 67
$ws.Range("H67").Value = 3236.2666
$ws.Range("I67").Value = 1666
$ws.Range("J67").Value = 3477.8462
$ws.Range("K67").Value = 4998
$ws.Range("L67").Value = 10433.5386
$ws.Range("M67").Value = -4062
$ws.Range("N67").Value = -12305.5386

# Row 68
$ws.Range("H68").Value = 600
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 200
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 600
$ws.Range("M68").Value = -2189
$ws.Range("N68").Value = -2222

# Row 70
$ws.Range("H70").Value = 3060
$ws.Range("I70").Value = 600
$ws.Range("J70").Value = 3675
$ws.Range("K70").Value = 1800
$ws.Range("L70").Value = 11025
$ws.Range("M70").Value = -1485
$ws.Range("N70").Value = -11655

# Row 71
$ws.Range("H71").Value = 600
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 200
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 1800
$ws.Range("M71").Value = -4944
$ws.Range("N71").Value = -9912

# Row 73
$ws.Range("H73").Value = 3060
$ws.Range("I73").Value = 600
$ws.Range("J73").Value = 3675
$ws.Range("K73").Value = 1800
$ws.Range("L73").Value = 11025
$ws.Range("M73").Value = -708
$ws.Range("N73").Value = -13209

# Row 121
$ws.Range("H121").Value = 1968.3405
$ws.Range("I121").Value = 3082.7273
$ws.Range("J121").Value = 1627.8334
$ws.Range("K121").Value = 9248.1819
$ws.Range("L121").Value = 4883.5002
$ws.Range("M121").Value = -7938.1819
$ws.Range("N121").Value = -7503.5002

# Row 131
$ws.Range("H131").Value = 1669.3
$ws.Range("I131").Value = 8226
$ws.Range("J131").Value = 1283.6118
$ws.Range("K131").Value = 24678
$ws.Range("L131").Value = 3850.8354
$ws.Range("M131").Value = -19638
$ws.Range("N131").Value = -13930.8354

$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 900.2222
$ws.Range("I22").Value = 933.3333
$ws.Range("J22").Value = 834
$ws.Range("K22").Value = 933.3333
$ws.Range("L22").Value = 834
$ws.Range("M22").Value = -638.3333
$ws.Range("N22").Value = -1424

# Row 27
$ws.Range("H27").Value = 900.2222
$ws.Range("I27").Value = 933.3333
$ws.Range("J27").Value = 834
$ws.Range("K27").Value = 933.3333
$ws.Range("L27").Value = 834
$ws.Range("M27").Value = -826.3333
$ws.Range("N27").Value = -1048

# Row 136
$ws.Range("H136").Value = 4095.4
$ws.Range("I136").Value = 2565.2415
$ws.Range("J136").Value = 5802.115
$ws.Range("K136").Value = 7695.7245
$ws.Range("L136").Value = 17406.345
$ws.Range("M136").Value = -5145.7245
$ws.Range("N136").Value = -22506.345

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 336.65216
$ws.Range("I113").Value = 312.78946
$ws.Range("K113").Value = 938.3683800000001
$ws.Range("M113").Value = 1231.63162

# Row 136
$ws.Range("H136").Value = 19917148
$ws.Range("I136").Value = 24415864
$ws.Range("J136").Value = 12823020
$ws.Range("K136").Value = 73247592
$ws.Range("L136").Value = 38469060
$ws.Range("M136").Value = -38474160
